$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.01365594787041438
$ws.Range("D2").Value = 104
$ws.Range("G2").Value = 189

# Row 5
$ws.Range("B5").Value = 0.006067347585136266
$ws.Range("C5").Value = 0.003836920719329858
$ws.Range("D5").Value = 249
$ws.Range("E5").Value = 220
$ws.Range("G5").Value = 270

# Row 7
$ws.Range("B7").Value = 0.00005722375326022583
$ws.Range("C7").Value = 0.00003851297747070992
$ws.Range("D7").Value = 855
$ws.Range("E7").Value = 784

# Row 8
$ws.Range("B8").Value = 0.09225411966263156
$ws.Range("C8").Value = 0.08077165916344395
$ws.Range("D8").Value = 145
$ws.Range("E8").Value = 132

# Row 9
$ws.Range("B9").Value = 0.0383531350879566
$ws.Range("C9").Value = 0.03226381414977151
$ws.Range("D9").Value = 286
$ws.Range("E9").Value = 255

# Row 10
$ws.Range("B10").Value = 0.3970927852420561
$ws.Range("C10").Value = 0.3763010113054139
$ws.Range("D10").Value = 42
$ws.Range("E10").Value = 36
$ws.Range("G10").Value = 124

# Row 11
$ws.Range("B11").Value = 0.000000002477806270654973
$ws.Range("C11").Value = 0.000000001101265678433875
$ws.Range("D11").Value = 1181
$ws.Range("E11").Value = 1098
$ws.Range("F11").Value = 1488
$ws.Range("G11").Value = 1481

# Row 13
$ws.Range("B13").Value = 0.507471706605126
$ws.Range("C13").Value = 0.5004452583181059
$ws.Range("D13").Value = 90
$ws.Range("E13").Value = 81

# Row 14
$ws.Range("C14").Value = 0.7589543307708074
$ws.Range("D14").Value = 32

# Row 16
$ws.Range("C16").Value = 0.6250617265957956
$ws.Range("D16").Value = 64
